# This script re-applies the betexplorer scraper's update for
# saudi-arabia_saudi-professional-league_2023-2024.xlsx:
#  - several existing match rows had their F:V (match details/odds) content
#    rotated/swapped between rows (the A:E "index/date" columns stay put)
#  - four brand new match rows (119-122) were appended at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($row) {
    # returns a 17-element array for columns F..V (index 0..16)
    $vals = @()
    for ($i = 0; $i -lt 17; $i++) {
        $col = 6 + $i  # F = 6
        $vals += $ws.Cells.Item($row, $col).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($i = 0; $i -lt 17; $i++) {
        $col = 6 + $i  # F = 6
        $ws.Cells.Item($row, $col).Value2 = $vals[$i]
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowValues $rowA
    $b = Get-RowValues $rowB
    Set-RowValues $rowA $b
    Set-RowValues $rowB $a
}

# --- 3-way rotation among rows 14, 15, 16 ---
# new14 = old15 ; new15 = old16 ; new16 = old14
$r14 = Get-RowValues 14
$r15 = Get-RowValues 15
$r16 = Get-RowValues 16
Set-RowValues 14 $r15
Set-RowValues 15 $r16
Set-RowValues 16 $r14

# --- simple pairwise swaps ---
Swap-Rows 81 82
Swap-Rows 83 84
Swap-Rows 85 86
Swap-Rows 98 99
Swap-Rows 105 106
Swap-Rows 110 111
Swap-Rows 116 117

# --- append 4 new rows (118-121 zero-based "Indice", sheet rows 119-122) ---
$newRows = @(
    @{ Idx = 118; Date = 45254.66666666666; F = "Al Taawon"; G = 1; H = "Al Riyadh"; I = 2;
       J = 1.35; K = "17/11/2023 16:42"; L = 1.37; M = "24/11/2023 15:50";
       N = 5.11; O = "17/11/2023 16:42"; P = 5.32; Q = "24/11/2023 15:54";
       R = 6.86; S = "17/11/2023 16:42"; T = 7.56; U = "24/11/2023 15:54";
       V = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taawon-al-riyadh/SbM0IOBM/" },
    @{ Idx = 119; Date = 45254.66666666666; F = "Al Ettifaq"; G = 1; H = "Al Ittihad"; I = 1;
       J = 3.84; K = "17/11/2023 16:42"; L = 3.08; M = "24/11/2023 15:50";
       N = 3.87; O = "17/11/2023 16:42"; P = 3.34; Q = "24/11/2023 15:50";
       R = 1.79; S = "17/11/2023 16:42"; T = 2.38; U = "24/11/2023 15:50";
       V = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ettifaq-fc-al-ittihad/KYIhKpsA/" },
    @{ Idx = 120; Date = 45254.79166666666; F = "Al Nassr"; G = 3; H = "Al Akhdoud"; I = 0;
       J = 1.12; K = "17/11/2023 19:43"; L = 1.24; M = "24/11/2023 18:55";
       N = 9.470000000000001; O = "17/11/2023 19:43"; P = 7.01; Q = "24/11/2023 18:57";
       R = 12.33; S = "17/11/2023 19:43"; T = 9.31; U = "24/11/2023 18:57";
       V = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-nassr-al-akhdoud/pbBvBt4d/" },
    @{ Idx = 121; Date = 45254.79166666666; F = "Al Taee"; G = 4; H = "Al Raed"; I = 3;
       J = 2.75; K = "17/11/2023 19:43"; L = 4.15; M = "24/11/2023 18:57";
       N = 3.67; O = "17/11/2023 19:43"; P = 3.63; Q = "24/11/2023 18:58";
       R = 2.3; S = "17/11/2023 19:43"; T = 1.88; U = "24/11/2023 18:53";
       V = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-raed/trKlLQR3/" }
)

$lastRow = 118
for ($n = 0; $n -lt $newRows.Count; $n++) {
    $targetRow = $lastRow + 1 + $n
    $srcRange = $ws.Range("A" + $lastRow + ":V" + $lastRow)
    $dstRange = $ws.Range("A" + $targetRow + ":V" + $targetRow)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0

    $d = $newRows[$n]
    $ws.Cells.Item($targetRow, 1).Value2 = $d.Idx            # A Indice
    $ws.Cells.Item($targetRow, 2).Value2 = "saudi-arabia"     # B pais
    $ws.Cells.Item($targetRow, 3).Value2 = "saudi-professional-league"  # C torneio
    $ws.Cells.Item($targetRow, 4).Value2 = "2023-2024"        # D temporada
    $ws.Cells.Item($targetRow, 5).Value2 = $d.Date            # E data_partida
    $ws.Cells.Item($targetRow, 6).Value2 = $d.F
    $ws.Cells.Item($targetRow, 7).Value2 = $d.G
    $ws.Cells.Item($targetRow, 8).Value2 = $d.H
    $ws.Cells.Item($targetRow, 9).Value2 = $d.I
    $ws.Cells.Item($targetRow, 10).Value2 = $d.J
    $ws.Cells.Item($targetRow, 11).Value2 = $d.K
    $ws.Cells.Item($targetRow, 12).Value2 = $d.L
    $ws.Cells.Item($targetRow, 13).Value2 = $d.M
    $ws.Cells.Item($targetRow, 14).Value2 = $d.N
    $ws.Cells.Item($targetRow, 15).Value2 = $d.O
    $ws.Cells.Item($targetRow, 16).Value2 = $d.P
    $ws.Cells.Item($targetRow, 17).Value2 = $d.Q
    $ws.Cells.Item($targetRow, 18).Value2 = $d.R
    $ws.Cells.Item($targetRow, 19).Value2 = $d.S
    $ws.Cells.Item($targetRow, 20).Value2 = $d.T
    $ws.Cells.Item($targetRow, 21).Value2 = $d.U
    $ws.Cells.Item($targetRow, 22).Value2 = $d.V
}
